# Update the "Förändrad" (changed) date column C for rows 2-23
# from 2023-10-25 (45224) to 2023-11-03 (45233).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 23; $row++) {
    $ws.Cells.Item($row, 3).Value = 45233
}
